$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new rows of training data (question -> intent) that were
# uploaded to the bottom of the sheet.
$rows = @(
    @('what''s the weather like', 'get_weather'),
    @('what''s the temperature', 'get_weather'),
    @('temperature', 'get_weather'),
    @('humidity', 'get_weather'),
    @('what is the humidity', 'get_weather'),
    @('how hot is it in Mumbai?', 'get_weather'),
    @('hello goodmorning', 'greeting'),
    @('hello jarvis how are you?', 'greeting'),
    @('good morning', 'greeting'),
    @('good afternoon', 'greeting'),
    @('how are you?', 'greeting')
)

$startRow = 191
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
}

# Match the author's final selection / scroll position.
$ws.Range("A186").Select() | Out-Null
